$d = $word.ActiveDocument

# 1) Replace the run text and the language tag on the paragraph mark + run
#    "Engenheiro(a) de Segurança do Trabalho" / pt-BR  ->  "FUNCAOHSE" / en-US
$d.Content.Find.Execute(
    "Engenheiro(a) de Segurança do Trabalho", $true, $false, $false, $false,
    $false, $true, 1, $false, "FUNCAOHSE", 2)

$d.Content.Find.Execute("FUNCAOHSE", $true, $false, $false, $false, $false, $true, 1, $false)
while ($d.Content.Find.Found) {
    $d.Content.Find.Parent.LanguageID = 1033
    $d.Content.Find.Execute("FUNCAOHSE", $true, $false, $false, $false, $false, $true, 1, $false)
}
